$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-339).
# All of these cells are being updated from serial 45186 (2023-09-17) to
# serial 45188 (2023-09-19). Use the last used row to stay safe even if the
# sheet's extents differ slightly from what we inspected.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 339 }

$newDate = [DateTime]::FromOADate(45188)

$ws.Range("C2:C$lastRow").Value = $newDate
